$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44176
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 4000
$ws.Range("P2").Value = 4000

$ws.Range("D3").Value = 44957
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 5000
$ws.Range("P3").Value = 5000

$ws.Range("D4").Value = 44312
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 4000
$ws.Range("P4").Value = 4000

$ws.Range("D5").Value = 44680
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 5000
$ws.Range("P5").Value = 5000

$ws.Range("D6").Value = 44749
$ws.Range("J6").Value = 65
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("P6").Value = 6000

$ws.Range("D7").Value = 44656
$ws.Range("J7").Value = 85
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = 5000
$ws.Range("P7").Value = 5000

$ws.Range("D8").Value = 44365
$ws.Range("J8").Value = 55
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = 5000
$ws.Range("P8").Value = 5000

$ws.Range("D9").Value = 44315
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 4000
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 4000
$ws.Range("P9").Value = 4000

$ws.Range("D10").Value = 44498
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 4000
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = 4000
$ws.Range("P10").Value = 4000

$ws.Range("D11").Value = 44679
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 5000
$ws.Range("P11").Value = 5000

$ws.Range("D12").Value = 44301
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("P12").Value = 3000

$ws.Range("D13").Value = 44504
$ws.Range("J13").Value = 55
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000

$ws.Range("D14").Value = 44390
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 6000
$ws.Range("P14").Value = 6000

$ws.Range("D15").Value = 44959
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = 5000
$ws.Range("P15").Value = 5000

$ws.Range("D16").Value = 44956
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = 5000
$ws.Range("P16").Value = 5000

$ws.Range("D17").Value = 44316
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = 4000
$ws.Range("P17").Value = 4000

$ws.Range("D18").Value = 44508
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 4000
$ws.Range("P18").Value = 4000

$ws.Range("D19").Value = 44966
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 5000
$ws.Range("P19").Value = 5000

$ws.Range("D20").Value = 44280
$ws.Range("J20").Value = 55
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = 4000
$ws.Range("P20").Value = 4000

$ws.Range("D21").Value = 44649
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 5000
$ws.Range("P21").Value = 5000

$ws.Range("D22").Value = 44509
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = 4000
$ws.Range("P22").Value = 4000

$ws.Range("D23").Value = 44259
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 4000
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = 4000
$ws.Range("P23").Value = 4000

$ws.Range("D24").Value = 44497
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 4000
$ws.Range("P24").Value = 4000

$ws.Range("D25").Value = 44781
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 5000
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 5000
$ws.Range("P25").Value = 5000

$ws.Range("D26").Value = 44291
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = 4000
$ws.Range("L26").Value = 4000
$ws.Range("M26").Value = 4000
$ws.Range("P26").Value = 4000

$ws.Range("D27").Value = 44777
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 5000
$ws.Range("P27").Value = 5000

$ws.Range("D28").Value = 44313
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("M28").Value = 4000
$ws.Range("P28").Value = 4000
